# Adds three rows of sample/possible values below the existing header row
# (name | border_type), giving the template concrete examples: a blank
# "name" column paired with border_type values "standoff", "ab", "pab".
#
# Note: assigning Value = "" is treated as "clear the cell" (it leaves no
# cell behind at all), so an explicit empty-string text cell in column A
# is produced via the classic leading-apostrophe text marker ("'"), which
# Excel stores as an empty string. Resetting Style back to "Normal"
# afterwards drops the quote-prefix formatting that the apostrophe trick
# applies, so column A ends up as a plain, unstyled empty text cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "standoff"

$ws.Range("A3").Value = "'"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "ab"

$ws.Range("A4").Value = "'"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "pab"
